$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Figure 2 legend paragraph: trim/rewrite the "x-axis" sentence
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "The x-axis shows counts of SSU sequences normalized to the average number of reads acquired for each sample.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The x-axis shows normalized counts of SSU sequences.", 2) | Out-Null

# "... belong to the same phylum or higher rank ..." -> drop "phylum or"
$d.Content.Find.Execute(
    "belong to the same phylum or higher rank",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "belong to the same higher rank", 2) | Out-Null

# "Taxa are numbered in the plot where they are highly abundant with
# corresponding numbered boxes in the legend." ->
# "Abundant Taxa are numbered in the plot with corresponding numbered
# boxes in the legend."
$d.Content.Find.Execute(
    "square bracket in the legend. Taxa are numbered in the plot where they are highly abundant with corresponding numbered boxes in the legend.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "square bracket in the legend. Abundant taxa are numbered in the plot with corresponding numbered boxes in the legend.", 2) | Out-Null

# Append the new "(C) Composition of abundant bacterial classes..."
# sentence at the end of the same paragraph. Use a placeholder around the
# "C" so it lands in the run produced by the Find/Replace (inheriting the
# paragraph's Times New Roman formatting), then bold just that character.
$d.Content.Find.Execute(
    "corresponding numbered boxes in the legend.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "corresponding numbered boxes in the legend. (#C#) Composition of abundant bacterial classes. SSU sequences were classified to the genus level or to the lowest rank with bootstap confidence >85% (see materials and methods).", 2) | Out-Null

$marker = $d.Content
$marker.Find.Execute("(#C#)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cOnly = $d.Range($marker.Start + 2, $marker.Start + 3)
$cOnly.Font.Bold = 1
$after = $d.Range($marker.Start + 3, $marker.Start + 4)
$after.Text = ""
$before = $d.Range($marker.Start + 1, $marker.Start + 2)
$before.Text = ""
